# Add an architecture-diagram of rounded-rectangle boxes, connectors,
# a donut "START" ring and a "START" label to slide 7 ("Algorithms").
#
# EMU -> point conversion: PowerPoint COM shape geometry is expressed in
# points (1 pt = 12700 EMU). All offsets/extents below are written as
# "<emu>/12700" so the resulting OOXML reproduces the exact EMU values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# msoShapeRoundedRectangle = 5
# msoConnectorStraight = 1
# msoShapeDonut = 18
# msoTextBox AddTextbox orientation 1 = msoTextOrientationHorizontal
# msoArrowheadTriangle = 2
# msoThemeColorAccent1 = 5, Accent3 = 7, Accent6 = 10, Text1 = 13

# ---------------------------------------------------------------------
# Shape 1 (id 4): "Afgeronde rechthoek 3" - GraphFrame
# ---------------------------------------------------------------------
$sh = $s.Shapes.AddShape(5, 4737782/12700, 2515904/12700, 2501462/12700, 1145628/12700)
$sh.Name = "Afgeronde rechthoek 3"
$sh.Fill.ForeColor.ObjectThemeColor = 5
$sh.Line.ForeColor.ObjectThemeColor = 5
$sh.TextFrame.VerticalAnchor = 3
$tr = $sh.TextFrame.TextRange
$tr.Text = "GraphFrame"
$tr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# Shape 2 (id 5): "Rechte verbindingslijn met pijl 5"
# ---------------------------------------------------------------------
$cn = $s.Shapes.AddConnector(1, (7375881+693682)/12700, 3109738/12700, 7375881/12700, 3109738/12700)
$cn.Name = "Rechte verbindingslijn met pijl 5"
$cn.HorizontalFlip = $true
$cn.Height = 0
$cn.Line.EndArrowheadStyle = 2
$cn.Line.ForeColor.ObjectThemeColor = 5

# ---------------------------------------------------------------------
# Shape 3 (id 6): "Afgeronde rechthoek 8" - SecondFrame
# ---------------------------------------------------------------------
$sh = $s.Shapes.AddShape(5, 8689675/12700, 2515904/12700, 2501462/12700, 1145628/12700)
$sh.Name = "Afgeronde rechthoek 8"
$sh.Fill.ForeColor.ObjectThemeColor = 7
$sh.Line.ForeColor.ObjectThemeColor = 7
$sh.TextFrame.VerticalAnchor = 3
$tr = $sh.TextFrame.TextRange
$tr.Text = "`rSecondFrame`ry - edges`rz - vertices`r"
$tr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# Shape 4 (id 7): "Rechte verbindingslijn met pijl 10"
# ---------------------------------------------------------------------
$cn = $s.Shapes.AddConnector(1, (3450266+1166648)/12700, 3088718/12700, 3450266/12700, 3088718/12700)
$cn.Name = "Rechte verbindingslijn met pijl 10"
$cn.HorizontalFlip = $true
$cn.Height = 0
$cn.Line.BeginArrowheadStyle = 2
$cn.Line.EndArrowheadStyle = 2
$cn.Line.ForeColor.ObjectThemeColor = 5

# ---------------------------------------------------------------------
# Shape 5 (id 8): "Afgeronde rechthoek 12" - Graph JComponent (GraphDisplay)
# ---------------------------------------------------------------------
$sh = $s.Shapes.AddShape(5, 864721/12700, 2515904/12700, 2501462/12700, 1145628/12700)
$sh.Name = "Afgeronde rechthoek 12"
$sh.Fill.ForeColor.ObjectThemeColor = 5
$sh.Line.ForeColor.ObjectThemeColor = 5
$sh.TextFrame.VerticalAnchor = 3
$tr = $sh.TextFrame.TextRange
$tr.Text = "Graph JComponent`r(GraphDisplay)"
$tr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# Shape 6 (id 9): "Rechte verbindingslijn met pijl 14"
# ---------------------------------------------------------------------
$cn = $s.Shapes.AddConnector(1, (6030554+1)/12700, 3777145/12700, 6030554/12700, (3777145+1019503)/12700)
$cn.Name = "Rechte verbindingslijn met pijl 14"
$cn.HorizontalFlip = $true
$cn.Left = 6030554/12700
$cn.Width = 1/12700
$cn.Line.BeginArrowheadStyle = 2
$cn.Line.EndArrowheadStyle = 2
$cn.Line.ForeColor.ObjectThemeColor = 5

# ---------------------------------------------------------------------
# Shape 7 (id 10): "Afgeronde rechthoek 16" - Mouse Click Listener
# ---------------------------------------------------------------------
$sh = $s.Shapes.AddShape(5, 817424/12700, 4985837/12700, 2501462/12700, 1145628/12700)
$sh.Name = "Afgeronde rechthoek 16"
$sh.Fill.ForeColor.ObjectThemeColor = 7
$sh.Line.ForeColor.ObjectThemeColor = 7
$sh.TextFrame.VerticalAnchor = 3
$tr = $sh.TextFrame.TextRange
$tr.Text = "Mouse Click Listener"
$tr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# Shape 8 (id 11): "Rechte verbindingslijn met pijl 17"
# ---------------------------------------------------------------------
$cn = $s.Shapes.AddConnector(1, 2052385/12700, 3777145/12700, 2052385/12700, (3777145+1019503)/12700)
$cn.Name = "Rechte verbindingslijn met pijl 17"
$cn.Width = 0
$cn.Line.BeginArrowheadStyle = 2
$cn.Line.EndArrowheadStyle = 2
$cn.Line.ForeColor.ObjectThemeColor = 5

# ---------------------------------------------------------------------
# Shape 9 (id 12): "Afgeronde rechthoek 18" - Color buttons
# ---------------------------------------------------------------------
$sh = $s.Shapes.AddShape(5, 5053093/12700, 4985837/12700, 1923393/12700, 1145628/12700)
$sh.Name = "Afgeronde rechthoek 18"
$sh.Fill.ForeColor.ObjectThemeColor = 10
$sh.Line.ForeColor.ObjectThemeColor = 10
$sh.TextFrame.VerticalAnchor = 3
$tr = $sh.TextFrame.TextRange
$tr.Text = "Color buttons"
$tr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# Shape 10 (id 13): "Afgeronde rechthoek 25" - Color Array
# ---------------------------------------------------------------------
$sh = $s.Shapes.AddShape(5, 7507256/12700, 4985837/12700, 1923393/12700, 1145628/12700)
$sh.Name = "Afgeronde rechthoek 25"
$sh.Fill.ForeColor.ObjectThemeColor = 10
$sh.Line.ForeColor.ObjectThemeColor = 10
$sh.TextFrame.VerticalAnchor = 3
$tr = $sh.TextFrame.TextRange
$tr.Text = "Color Array"
$tr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------
# Shape 11 (id 14): "Rechte verbindingslijn met pijl 26"
# ---------------------------------------------------------------------
$cn = $s.Shapes.AddConnector(1, 7061878/12700, 5558651/12700, (7061878+396769)/12700, 5558651/12700)
$cn.Name = "Rechte verbindingslijn met pijl 26"
$cn.Height = 0
$cn.Line.BeginArrowheadStyle = 2
$cn.Line.EndArrowheadStyle = 2
$cn.Line.ForeColor.ObjectThemeColor = 5

# ---------------------------------------------------------------------
# Shape 12 (id 15): "Ring 29" (donut / START ring)
# ---------------------------------------------------------------------
$sh = $s.Shapes.AddShape(18, 8274516/12700, 1853754/12700, 3468414/12700, 2511967/12700)
$sh.Name = "Ring 29"
$sh.Adjustments.Item(1) = 0.04545
$sh.Fill.ForeColor.RGB = 255
$sh.Line.ForeColor.RGB = 192
$sh.TextFrame.VerticalAnchor = 3
$tr = $sh.TextFrame.TextRange
$tr.Text = ""
$tr.ParagraphFormat.Alignment = 2
$tr.Font.Color.ObjectThemeColor = 13

# ---------------------------------------------------------------------
# Shape 13 (id 16): "Tekstvak 31" - START label textbox
# ---------------------------------------------------------------------
$tb = $s.Shapes.AddTextbox(1, 9604074/12700, 3786924/12700, 809297/12700, 369332/12700)
$tb.Name = "Tekstvak 31"
$tb.TextFrame.TextRange.Text = "START"
$tb.TextFrame.TextRange.Font.Bold = $true
$tb.TextFrame.WordWrap = -1
$tb.TextFrame.AutoSize = 1
$tb.Fill.Visible = $false
$tb.Height = 369332/12700
